$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("B2", 14.25885541126553),
    @("C2", 8.460474104635257),
    @("E2", 12.36247422199909),
    @("F2", 16.86991607391233),
    @("G2", 52.2700144146375),
    @("H2", 20.16600593833469),
    @("K2", 10.62449821130868),
    @("L2", 10.07342784662966),
    @("M2", 15.15622141892267),
    @("B3", 14.12455322514939),
    @("C3", 8.435599950674264),
    @("E3", 12.38106974861166),
    @("F3", 15.89584955866808),
    @("G3", 52.24337459732921),
    @("H3", 20.19972502377299),
    @("K3", 10.53425400819322),
    @("L3", 10.08299292575401),
    @("M3", 15.14932502427314),
    @("B4", 14.04500536636511),
    @("C4", 8.419871745827367),
    @("E4", 12.39391830670356),
    @("F4", 15.26997757108491),
    @("G4", 52.23804067308948),
    @("H4", 20.22325666672233),
    @("K4", 10.48100879622683),
    @("L4", 10.09022174160067),
    @("M4", 15.14762361622899),
    @("B5", 14.01335714714178),
    @("C5", 8.413348297365529),
    @("E5", 12.39951424020416),
    @("F5", 15.00819731993403),
    @("G5", 52.23863685384855),
    @("H5", 20.23355624261685),
    @("K5", 10.45987805296979),
    @("L5", 10.09350883900024),
    @("M5", 15.14756936147675),
    @("B6", 14.00814938487088),
    @("C6", 8.412258206352817),
    @("E6", 12.40046519315299),
    @("F6", 14.96433081551593),
    @("G6", 52.23890304496864),
    @("H6", 20.23530934900252),
    @("K6", 10.45640420361034),
    @("L6", 10.09407528204262),
    @("M6", 15.14759900829673),
    @("B7", 14.04457539146961),
    @("C7", 8.41978422963729),
    @("E7", 12.39399231730846),
    @("F7", 15.26647399323137),
    @("G7", 52.23803750264869),
    @("H7", 20.22339269591947),
    @("K7", 10.48072149423546),
    @("L7", 10.09026469020898),
    @("M7", 15.14762029430712),
    @("B8", 14.211964833783),
    @("C8", 8.451992262267591),
    @("E8", 12.36858914578444),
    @("F8", 16.53996406344768),
    @("G8", 52.25854044822923),
    @("H8", 20.17704481955083),
    @("K8", 10.59294656424141),
    @("L8", 10.07644471931531),
    @("M8", 15.15331918187296),
    @("B9", 14.56171755672871),
    @("C9", 8.511521851754916),
    @("E9", 12.33011542407194),
    @("F9", 19.00274580682531),
    @("G9", 52.3862326502081),
    @("H9", 20.10863342985334),
    @("K9", 10.82915511710816),
    @("L9", 10.06008429310215),
    @("M9", 15.18448840870643),
    @("B10", 14.82954332540609),
    @("C10", 8.553024924660162),
    @("E10", 12.30874759732449),
    @("F10", 20.67494806633232),
    @("G10", 52.53331525099541),
    @("H10", 20.07212340455551),
    @("K10", 11.0110807017132),
    @("L10", 10.05458674273342),
    @("M10", 15.21942134550954),
    @("B11", 14.95325060322575),
    @("C11", 8.57141668312129),
    @("E11", 12.30052108335407),
    @("F11", 21.3917225636224),
    @("G11", 52.61174180614637),
    @("H11", 20.05850926122921),
    @("K11", 11.09534256090918),
    @("L11", 10.05349550220414),
    @("M11", 15.23788676977002),
    @("B12", 15.00032282350773),
    @("C12", 8.578310549297917),
    @("E12", 12.29762035951646),
    @("F12", 21.65686569030329),
    @("G12", 52.6430884537471),
    @("H12", 20.05378512444944),
    @("K12", 11.12743905124246),
    @("L12", 10.05328429753815),
    @("M12", 15.24524545455976),
    @("B13", 14.9901755561108),
    @("C13", 8.5768289837733),
    @("E13", 12.29823554890094),
    @("F13", 21.60004134736742),
    @("G13", 52.63626423578954),
    @("H13", 20.05478336283447),
    @("K13", 11.12051856387794),
    @("L13", 10.05332081066076),
    @("M13", 15.2436444080085),
    @("B14", 14.95711895614968),
    @("C14", 8.571985265246985),
    @("E14", 12.30027814259238),
    @("F14", 21.4136618050453),
    @("G14", 52.61428773387117),
    @("H14", 20.05811195751433),
    @("K14", 11.09797954302827),
    @("L14", 10.0534740815196),
    @("M14", 15.23848485708628),
    @("B15", 14.93689913867748),
    @("C15", 8.569009126875779),
    @("E15", 12.30155721129977),
    @("F15", 21.29868154950795),
    @("G15", 52.60104086680391),
    @("H15", 20.06020699613464),
    @("K15", 11.08419744506818),
    @("L15", 10.05359425262815),
    @("M15", 15.235372054789),
    @("B16", 14.82149282243751),
    @("C16", 8.551813114122453),
    @("E16", 12.30931525121582),
    @("F16", 20.62722412089977),
    @("G16", 52.52842101624113),
    @("H16", 20.07307343409694),
    @("K16", 11.00560186354006),
    @("L16", 10.05468636670201),
    @("M16", 15.21826605686832),
    @("B17", 14.75114420036262),
    @("C17", 8.541138527428373),
    @("E17", 12.31445693428863),
    @("F17", 20.20408069597325),
    @("G17", 52.48681598048213),
    @("H17", 20.08173401580996),
    @("K17", 10.95775125494419),
    @("L17", 10.05571692082847),
    @("M17", 15.20842869916594),
    @("B18", 14.71086026833818),
    @("C18", 8.53495286523747),
    @("E18", 12.31755492127958),
    @("F18", 19.95656407809801),
    @("G18", 52.46397058932131),
    @("H18", 20.0869971721667),
    @("K18", 10.93037199502272),
    @("L18", 10.05644239346364),
    @("M18", 15.20301306782179),
    @("B19", 14.69725280599157),
    @("C19", 8.532850634131723),
    @("E19", 12.31862800882215),
    @("F19", 19.87204792380568),
    @("G19", 52.45642204957537),
    @("H19", 20.08882756801017),
    @("K19", 10.9211272774956),
    @("L19", 10.0567108412475),
    @("M19", 15.20122120269079),
    @("B20", 14.75861473893146),
    @("C20", 8.5422796117654),
    @("E20", 12.31389504129609),
    @("F20", 20.24955283636154),
    @("G20", 52.49113269487336),
    @("H20", 20.0807829099123),
    @("K20", 10.96283041983255),
    @("L20", 10.05559348472077),
    @("M20", 15.20945082660107),
    @("B21", 14.96682265744757),
    @("C21", 8.573409905293948),
    @("E21", 12.29967236521531),
    @("F21", 21.46857628470577),
    @("G21", 52.62069810666108),
    @("H21", 20.0571225600452),
    @("K21", 11.10459490976496),
    @("L21", 10.05342358500335),
    @("M21", 15.23999043584136),
    @("B22", 15.10420197486087),
    @("C22", 8.593343070026076),
    @("E22", 12.29162699371624),
    @("F22", 22.22866616901552),
    @("G22", 52.71497858331955),
    @("H22", 20.04417300701227),
    @("K22", 11.19833115063359),
    @("L22", 10.05318262092738),
    @("M22", 15.26208245899186),
    @("B23", 15.03077470661146),
    @("C23", 8.582742250023028),
    @("E23", 12.29580670477447),
    @("F23", 21.82633154458858),
    @("G23", 52.66378388102771),
    @("H23", 20.05085422000939),
    @("K23", 11.14821227358196),
    @("L23", 10.05320375631085),
    @("M23", 15.25009777647213),
    @("B24", 14.75523680549799),
    @("C24", 8.541763879212461),
    @("E24", 12.31414863089291),
    @("F24", 20.22900810905287),
    @("G24", 52.48917776246896),
    @("H24", 20.08121202001284),
    @("K24", 10.9605337201649),
    @("L24", 10.05564887587393),
    @("M24", 15.2089879749161),
    @("B25", 14.46502682916583),
    @("C25", 8.495808258290079),
    @("E25", 12.33931085519139),
    @("F25", 18.34778573295695),
    @("G25", 52.34232214264208),
    @("H25", 20.12472881934566),
    @("K25", 10.76367452601929),
    @("L25", 10.06336273488915),
    @("M25", 15.17393071144759)
)

foreach ($change in $changes) {
    $ws.Range($change[0]).Value = $change[1]
}

Write-Host "Applied" $changes.Count "cell updates"
